# Scheduled-runner price refresh: overwrite cached market-price/profit
# columns (H:N) per row across all 8 Leve-profit sheets, per the source diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7953.75
$ws.Range("J40").Value = 7275.6665
$ws.Range("L40").Value = 7275.6665
$ws.Range("N40").Value = -7625.6665
$ws.Range("H106").Value = 9999.5
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H111").Value = 2706.75
$ws.Range("J111").Value = 7499
$ws.Range("L111").Value = 22497
$ws.Range("N111").Value = -28631
$ws.Range("H132").Value = 321881.56
$ws.Range("I132").Value = 362505.3
$ws.Range("K132").Value = 1087515.9
$ws.Range("M132").Value = -1084985.9
$ws.Range("H137").Value = 12636.263
$ws.Range("I137").Value = 13903.6
$ws.Range("J137").Value = 11228.111
$ws.Range("K137").Value = 41710.8
$ws.Range("L137").Value = 33684.333
$ws.Range("M137").Value = -39160.8
$ws.Range("N137").Value = -38784.333
$ws.Range("H141").Value = 2093.8
$ws.Range("I141").Value = 2093.8
$ws.Range("K141").Value = 6281.400000000001
$ws.Range("M141").Value = -1101.400000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4633979
$ws.Range("I32").Value = 2888.842
$ws.Range("K32").Value = 2888.842
$ws.Range("M32").Value = -2601.842
$ws.Range("H36").Value = 14525
$ws.Range("I36").Value = 9700
$ws.Range("K36").Value = 9700
$ws.Range("M36").Value = -9354
$ws.Range("H58").Value = 8000
$ws.Range("J58").Value = 8000
$ws.Range("L58").Value = 8000
$ws.Range("N58").Value = -8860
$ws.Range("H61").Value = 4392.5557
$ws.Range("I61").Value = 3785
$ws.Range("K61").Value = 3785
$ws.Range("M61").Value = -3573
$ws.Range("H74").Value = 4676.381
$ws.Range("I74").Value = 4356.5
$ws.Range("K74").Value = 4356.5
$ws.Range("M74").Value = -3482.5
$ws.Range("H77").Value = 4676.381
$ws.Range("I77").Value = 4356.5
$ws.Range("K77").Value = 21782.5
$ws.Range("M77").Value = -17414.5
$ws.Range("H132").Value = 772763.3
$ws.Range("I132").Value = 843765.3
$ws.Range("K132").Value = 2531295.9
$ws.Range("M132").Value = -2528765.9
$ws.Range("H136").Value = 4392.5557
$ws.Range("I136").Value = 3785
$ws.Range("K136").Value = 11355
$ws.Range("M136").Value = -8805

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3560.6667
$ws.Range("I22").Value = 272.8
$ws.Range("J22").Value = 20000
$ws.Range("K22").Value = 272.8
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = -99.80000000000001
$ws.Range("N22").Value = -20346
$ws.Range("H80").Value = 12831228
$ws.Range("J80").Value = 18533332
$ws.Range("L80").Value = 18533332
$ws.Range("N80").Value = -18535328
$ws.Range("H83").Value = 12831228
$ws.Range("J83").Value = 18533332
$ws.Range("L83").Value = 92666660
$ws.Range("N83").Value = -92676644
$ws.Range("H105").Value = 4162.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5424
$ws.Range("I31").Value = 1318.375
$ws.Range("J31").Value = 8708.5
$ws.Range("K31").Value = 1318.375
$ws.Range("L31").Value = 8708.5
$ws.Range("M31").Value = -1023.375
$ws.Range("N31").Value = -9298.5
$ws.Range("H34").Value = 5424
$ws.Range("I34").Value = 1318.375
$ws.Range("J34").Value = 8708.5
$ws.Range("K34").Value = 1318.375
$ws.Range("L34").Value = 8708.5
$ws.Range("M34").Value = -1116.375
$ws.Range("N34").Value = -9112.5
$ws.Range("H58").Value = 28583232
$ws.Range("I58").Value = 41673930
$ws.Range("K58").Value = 41673930
$ws.Range("M58").Value = -41673727
$ws.Range("H86").Value = 11788.637
$ws.Range("I86").Value = 7615.143
$ws.Range("K86").Value = 7615.143
$ws.Range("M86").Value = -6492.143
$ws.Range("H89").Value = 11788.637
$ws.Range("I89").Value = 7615.143
$ws.Range("K89").Value = 38075.715
$ws.Range("M89").Value = -32459.715
$ws.Range("H136").Value = 28583232
$ws.Range("I136").Value = 41673930
$ws.Range("K136").Value = 125021790
$ws.Range("M136").Value = -125019240

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 3799
$ws.Range("I8").Value = 3799
$ws.Range("K8").Value = 11397
$ws.Range("M8").Value = -11258
$ws.Range("H47").Value = 2044
$ws.Range("I47").Value = 609
$ws.Range("J47").Value = 3000.6667
$ws.Range("K47").Value = 1827
$ws.Range("L47").Value = 9002.000100000001
$ws.Range("M47").Value = -1396
$ws.Range("N47").Value = -9864.000100000001
$ws.Range("H116").Value = 3791.4167
$ws.Range("I116").Value = 3450.7
$ws.Range("K116").Value = 10352.1
$ws.Range("M116").Value = -6910.099999999999
$ws.Range("H119").Value = 3466.3333
$ws.Range("I119").Value = 3466.3333
$ws.Range("K119").Value = 10398.9999
$ws.Range("M119").Value = -5560.999899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 17853.285
$ws.Range("J15").Value = 17853.285
$ws.Range("L15").Value = 17853.285
$ws.Range("N15").Value = -18429.285
$ws.Range("H41").Value = 11790.8
$ws.Range("I41").Value = 16333
$ws.Range("J41").Value = 4977.5
$ws.Range("K41").Value = 16333
$ws.Range("L41").Value = 4977.5
$ws.Range("M41").Value = -15978
$ws.Range("N41").Value = -5687.5
$ws.Range("H54").Value = 12789.667
$ws.Range("J54").Value = 12789.667
$ws.Range("L54").Value = 12789.667
$ws.Range("N54").Value = -13569.667
$ws.Range("H81").Value = 17853.285
$ws.Range("J81").Value = 17853.285
$ws.Range("L81").Value = 17853.285
$ws.Range("N81").Value = -19849.285
$ws.Range("H84").Value = 17853.285
$ws.Range("J84").Value = 17853.285
$ws.Range("L84").Value = 53559.855
$ws.Range("N84").Value = -63543.855
$ws.Range("H132").Value = 23812666
$ws.Range("I132").Value = 28574084
$ws.Range("K132").Value = 85722252
$ws.Range("M132").Value = -85719722
$ws.Range("H138").Value = 75000
$ws.Range("J138").Value = 75000
$ws.Range("L138").Value = 75000
$ws.Range("N138").Value = -85280

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7302.72
$ws.Range("I61").Value = 6387.467
$ws.Range("K61").Value = 6387.467
$ws.Range("M61").Value = -6185.467
$ws.Range("H113").Value = 7302.72
$ws.Range("I113").Value = 6387.467
$ws.Range("K113").Value = 6387.467
$ws.Range("M113").Value = -4217.467
$ws.Range("H122").Value = 5555.857
$ws.Range("I122").Value = 4778.2
$ws.Range("K122").Value = 14334.6
$ws.Range("M122").Value = -11884.6
$ws.Range("H132").Value = 7499.25
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 7499.25
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 22497.75
$ws.Range("N132").Value = -27557.75
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 37045030
$ws.Range("J136").Value = 7661.8667
$ws.Range("L136").Value = 22985.6001
$ws.Range("N136").Value = -28085.6001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1663.4117
$ws.Range("I81").Value = 1481.5
$ws.Range("J81").Value = 2100
$ws.Range("K81").Value = 2963
$ws.Range("L81").Value = 4200
$ws.Range("M81").Value = -1902
$ws.Range("N81").Value = -6322
$ws.Range("H84").Value = 1663.4117
$ws.Range("I84").Value = 1481.5
$ws.Range("J84").Value = 2100
$ws.Range("K84").Value = 14815
$ws.Range("L84").Value = 21000
$ws.Range("M84").Value = -9511
$ws.Range("N84").Value = -31608
$ws.Range("H107").Value = 2256.1785
$ws.Range("I107").Value = 2016.4736
$ws.Range("J107").Value = 2762.2222
$ws.Range("K107").Value = 6049.4208
$ws.Range("L107").Value = 8286.6666
$ws.Range("M107").Value = -4129.4208
$ws.Range("N107").Value = -12126.6666
$ws.Range("H136").Value = 17253884
$ws.Range("I136").Value = 21749500
$ws.Range("K136").Value = 65248500
$ws.Range("M136").Value = -65245950
